$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "State" verbs (Present Perfect Continuous additions) into column C,
# rows 19-21, matching column E's style (text already present in E19:E21).
$ws.Range("C19").Value = " stay (останавливаться)"
$ws.Range("C20").Value = " exist (существовать)"
$ws.Range("C21").Value = " remain (оставаться)"

# Update selection to match the new active cell.
$ws.Range("I21").Select()
